# Scheduled runner update: refresh profit-calculation columns (H:N) on the
# Leve profit sheets (currentAveragePrice*, LevePrice*, LeveProfit*) with the
# latest market-board figures. A few rows also gain/lose an M (LeveProfitNQ)
# or N (LeveProfitHQ) cell depending on whether that recipe variant exists.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1546.5416
$ws.Range("I28").Value = 1133.3334
$ws.Range("J28").Value = 1959.75
$ws.Range("K28").Value = 1133.3334
$ws.Range("L28").Value = 1959.75
$ws.Range("M28").Value = -648.3334
$ws.Range("N28").Value = -2929.75

$ws.Range("H51").Value = 6057.8335
$ws.Range("I51").Value = 2800
$ws.Range("J51").Value = 6709.4
$ws.Range("K51").Value = 2800
$ws.Range("L51").Value = 6709.4
$ws.Range("M51").Value = -2316
$ws.Range("N51").Value = -7677.4

$ws.Range("H69").Value = 20166.334
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 20166.334
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 60499.00199999999
$ws.Range("M69").ClearContents()  # was -2445.5
$ws.Range("N69").Value = -62247.00199999999

$ws.Range("H72").Value = 20166.334
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 20166.334
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 181497.006
$ws.Range("M72").ClearContents()  # was -5590.5
$ws.Range("N72").Value = -190233.006

$ws.Range("H138").Value = 2629.0322
$ws.Range("I138").Value = 2093.5557
$ws.Range("J138").Value = 3370.4614
$ws.Range("K138").Value = 6280.6671
$ws.Range("L138").Value = 10111.3842
$ws.Range("M138").Value = -1140.6671
$ws.Range("N138").Value = -20391.3842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 114.31579
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 119.42857
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 119.42857
$ws.Range("M5").Value = 12
$ws.Range("N5").Value = -343.42857

$ws.Range("H122").Value = 1703.5428
$ws.Range("I122").Value = 1595.6666
$ws.Range("J122").Value = 1938.909
$ws.Range("K122").Value = 4786.9998
$ws.Range("L122").Value = 5816.727000000001
$ws.Range("M122").Value = -2336.9998
$ws.Range("N122").Value = -10716.727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 114.31579
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 119.42857
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 119.42857
$ws.Range("M4").Value = 15
$ws.Range("N4").Value = -349.42857

$ws.Range("H134").Value = 772456.5600000001
$ws.Range("I134").Value = 854011.2
$ws.Range("J134").Value = 5842.8
$ws.Range("K134").Value = 2562033.6
$ws.Range("L134").Value = 17528.4
$ws.Range("M134").Value = -2559498.6
$ws.Range("N134").Value = -22598.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 27874
$ws.Range("J63").Value = 27874
$ws.Range("L63").Value = 27874
$ws.Range("N63").Value = -29246

$ws.Range("H66").Value = 27874
$ws.Range("J66").Value = 27874
$ws.Range("L66").Value = 83622
$ws.Range("N66").Value = -90486

$ws.Range("H129").Value = 35299.875
$ws.Range("J129").Value = 35299.875
$ws.Range("L129").Value = 35299.875
$ws.Range("N129").Value = -45299.875

$ws.Range("H134").Value = 1414.8667
$ws.Range("I134").Value = 1444.5
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 4333.5
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -1798.5
$ws.Range("N134").Value = -8070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 4950.6665
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()  # was -5225

$ws.Range("H72").Value = 4950.6665
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()  # was -14052

$ws.Range("H74").Value = 8012.6665
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 8012.6665
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 24037.9995
$ws.Range("M74").ClearContents()  # was 761
$ws.Range("N74").Value = -26159.9995

$ws.Range("H77").Value = 8012.6665
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 8012.6665
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 72113.9985
$ws.Range("M77").ClearContents()  # was 4404
$ws.Range("N77").Value = -82721.9985

$ws.Range("H80").Value = 3125.7856
$ws.Range("I80").Value = 800
$ws.Range("K80").Value = 2400
$ws.Range("M80").Value = -1464

$ws.Range("H83").Value = 3125.7856
$ws.Range("I83").Value = 800
$ws.Range("K83").Value = 7200
$ws.Range("M83").Value = -2520

$ws.Range("H131").Value = 1086.7037
$ws.Range("I131").Value = 775.25
$ws.Range("J131").Value = 1102.883
$ws.Range("K131").Value = 2325.75
$ws.Range("L131").Value = 3308.649
$ws.Range("M131").Value = 2714.25
$ws.Range("N131").Value = -13388.649

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 16955.666
$ws.Range("I69").Value = 10182
$ws.Range("J69").Value = 18310.4
$ws.Range("K69").Value = 10182
$ws.Range("L69").Value = 18310.4
$ws.Range("M69").Value = -9433
$ws.Range("N69").Value = -19808.4

$ws.Range("H70").Value = 4594.8286
$ws.Range("I70").Value = 4490.913
$ws.Range("J70").Value = 4794
$ws.Range("K70").Value = 4490.913
$ws.Range("L70").Value = 4794
$ws.Range("M70").Value = -4220.913
$ws.Range("N70").Value = -5334

$ws.Range("H72").Value = 16955.666
$ws.Range("I72").Value = 10182
$ws.Range("J72").Value = 18310.4
$ws.Range("K72").Value = 30546
$ws.Range("L72").Value = 54931.2
$ws.Range("M72").Value = -26802
$ws.Range("N72").Value = -62419.2

$ws.Range("H73").Value = 4594.8286
$ws.Range("I73").Value = 4490.913
$ws.Range("J73").Value = 4794
$ws.Range("K73").Value = 4490.913
$ws.Range("L73").Value = 4794
$ws.Range("M73").Value = -3554.913
$ws.Range("N73").Value = -6666

$ws.Range("H113").Value = 1698
$ws.Range("I113").Value = 1400
$ws.Range("J113").Value = 1896.6666
$ws.Range("K113").Value = 1400
$ws.Range("L113").Value = 1896.6666
$ws.Range("M113").Value = 770
$ws.Range("N113").Value = -6236.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 9568.538
$ws.Range("I68").Value = 26375
$ws.Range("J68").Value = 2099
$ws.Range("K68").Value = 26375
$ws.Range("L68").Value = 2099
$ws.Range("M68").Value = -25626
$ws.Range("N68").Value = -3597

$ws.Range("H71").Value = 9568.538
$ws.Range("I71").Value = 26375
$ws.Range("J71").Value = 2099
$ws.Range("K71").Value = 131875
$ws.Range("L71").Value = 10495
$ws.Range("M71").Value = -128131
$ws.Range("N71").Value = -17983

$ws.Range("H132").Value = 3892.9
$ws.Range("I132").Value = 4029.8333
$ws.Range("J132").Value = 3482.1
$ws.Range("K132").Value = 12089.4999
$ws.Range("L132").Value = 10446.3
$ws.Range("M132").Value = -9559.499899999999
$ws.Range("N132").Value = -15506.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5433.3335
$ws.Range("I62").Value = 6500
$ws.Range("J62").Value = 3300
$ws.Range("K62").Value = 6500
$ws.Range("L62").Value = 3300
$ws.Range("M62").Value = -5876
$ws.Range("N62").Value = -4548

$ws.Range("H65").Value = 5433.3335
$ws.Range("I65").Value = 6500
$ws.Range("J65").Value = 3300
$ws.Range("K65").Value = 32500
$ws.Range("L65").Value = 16500
$ws.Range("M65").Value = -29380
$ws.Range("N65").Value = -22740

$ws.Range("H136").Value = 8794
$ws.Range("I136").Value = 6878.6665
$ws.Range("J136").Value = 10120
$ws.Range("K136").Value = 20635.9995
$ws.Range("L136").Value = 30360
$ws.Range("M136").Value = -18085.9995
$ws.Range("N136").Value = -35460
